# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 2
    6  = 2
    7  = 1
    8  = 0
    9  = 1
    10 = 0
    11 = 2
    12 = 2
    13 = 2
    14 = 0
    15 = 6
    16 = 4
    17 = 2
    18 = 1
    19 = 1
    20 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
